# Adds a new forecast-vintage column Z (dated 2020-04-27) and a new
# observation row 38 (dated 2020-05-11) to both the "cases" and "deaths"
# sheets, and backfills the previously-empty "Observed" value for
# 2020-04-27 (row 24, column B).
#
# Column Z holds the forecasts made "as of" 2020-04-27 for every future
# date; following the existing diagonal pattern in the sheet, the first
# value lands in row 25 (2020-04-28) and the new row 38 (2020-05-11) is
# only populated in column Z.

$wb = $excel.ActiveWorkbook

# Z-column ("as of 2020-04-27") values per sheet, keyed by row number.
$casesZ = @{
    25 = 71283
    26 = 77145
    27 = 84760
    28 = 91961
    29 = 100507
    30 = 107551
    31 = 118208
    32 = 127747
    33 = 136769
    34 = 145151
    35 = 152631
    36 = 158104
    37 = 165041
    38 = 171322
}

$deathsZ = @{
    25 = 4873
    26 = 5224
    27 = 5679
    28 = 6071
    29 = 6559
    30 = 6870
    31 = 7445
    32 = 7870
    33 = 8261
    34 = 8643
    35 = 8981
    36 = 9232
    37 = 9544
    38 = 9823
}

$sheetSpecs = @(
    @{ Name = "cases";  ObservedB24 = 66501; ZValues = $casesZ },
    @{ Name = "deaths"; ObservedB24 = 4543;  ZValues = $deathsZ }
)

foreach ($spec in $sheetSpecs) {
    $ws = $wb.Worksheets.Item($spec.Name)

    # --- Header: Z1 = "2020-04-27" (plain text, not an Excel date) ---
    $z1 = $ws.Cells.Item(1, 26)
    $z1.NumberFormat = "@"
    $z1.Value = "2020-04-27"
    $z1.ClearFormats()

    # --- Blank out Z2:Z24 explicitly so those cells exist (empty) ---
    $zBlankRange = $ws.Range($ws.Cells.Item(2, 26), $ws.Cells.Item(24, 26))
    $zBlankRange.NumberFormat = "General"
    $zBlankRange.ClearFormats()

    # --- Backfill the observed value for 2020-04-27 (row 24, col B) ---
    $ws.Cells.Item(24, 2).Value = $spec.ObservedB24

    # --- New forecast-vintage values in column Z, rows 25-38 ---
    foreach ($r in ($spec.ZValues.Keys | Sort-Object)) {
        $ws.Cells.Item($r, 26).Value = $spec.ZValues[$r]
    }

    # --- New row 38: A38 = "2020-05-11" (text), B38:Y38 blank, Z38 set above ---
    $a38 = $ws.Cells.Item(38, 1)
    $a38.NumberFormat = "@"
    $a38.Value = "2020-05-11"
    $a38.ClearFormats()

    $row38BlankRange = $ws.Range($ws.Cells.Item(38, 2), $ws.Cells.Item(38, 25))
    $row38BlankRange.NumberFormat = "General"
    $row38BlankRange.ClearFormats()
}
